$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '26.159.30'
$c.Style = 'Normal'
$ws.Range('E2').Value = '  -1.23%  '
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '1.657.82'
$c.Style = 'Normal'
$ws.Range('E3').Value = '  -1.22%  '
$ws.Range('E4').Value = '  +0.35%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '215.95'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -1.73%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '0.5197'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -2.71%  '
$ws.Range('E7').Value = '  +0.34%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.2623'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  -3.31%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.06259'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  -2.45%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '20.71'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  -5.72%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.07719'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  -0.99%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '1.668.07'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  -0.73%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '4.413'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  -2.30%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '1.882.97'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  -1.38%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '0.5414'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  -3.36%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '0.0₅8114'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  -2.88%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '64.32'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  -2.26%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '26.174.06'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  -1.34%  '
$ws.Range('E19').Value = '  +0.42%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '4.620'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -3.99%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '191.37'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  -1.18%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '10.03'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -2.82%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '6.054'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  -4.30%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '1.007'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  +0.49%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '139.79'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  -1.78%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '0.1223'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  -4.87%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '7.171'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -3.42%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '16.05'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  -1.56%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '1.413'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  -2.10%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '0.05983'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  -5.12%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '1.270'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  -1.36%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '3.554'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -1.54%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '3.232'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  -6.70%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '1.609'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  -5.25%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '0.9644'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  -4.79%  '
$ws.Range('E36').Value = '  -0.13%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '2.770'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  -0.62%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '0.5659'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  -8.03%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '6.008'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -2.36%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '0.01596'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  -2.38%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.8551'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  -1.37%  '
$ws.Range('E42').Value = '  +0.40%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '1.011.95'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  -7.38%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '100.14'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  -0.49%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '1.798.28'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -1.42%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '0.0₈107'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -5.15%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '56.63'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  -3.68%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '1.006'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +0.36%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '7.927'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  -3.38%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '0.05173'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  -0.72%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '0.4204'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -0.75%  '
